$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- CRM-3026: Add a "Current Status" column ahead of the existing
# internal-status column, and relabel that existing column "Internal Status".
# Everything else shifts one column to the right. ---

# Insert a new column at A; old A..Q become B..R.
$ws.Columns("A:A").Insert()

# Relabel the (now shifted) old "Status" header -> "Internal Status".
# (Its value cell below still carries {bookings:internal_status} untouched.)
$ws.Range("B1").Value = "Internal Status"

# New first column: header + merge-field value.
$ws.Range("A1").Value = "Current Status"
$ws.Range("A2").Value = "{bookings:current_status}"

# Header row formatting to match the existing bold headers.
$ws.Range("A1:B1").Font.Bold = $true

# Data row: "Internal Status" value cell is center-aligned like the other
# status-style value cells; "Current Status" value cell stays general/left.
$ws.Range("B2").HorizontalAlignment = -4108

# Column widths (new/changed columns only - the rest shift along with
# the column insert above and keep their original widths).
$ws.Range("A1").EntireColumn.ColumnWidth = 21.83
$ws.Range("B1").EntireColumn.ColumnWidth = 27.39
$ws.Range("C1").EntireColumn.ColumnWidth = 20.83

# Restore the view to the top-left with A2 selected (no more scrolled-right
# / stale P4 selection from before the column was inserted).
$ws.Range("A2").Select()
